# Apply the edits described by the diff to vocabulary.xlsx
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# 1. A1 header text change
$ws.Range("A1").Value = "ConceptScheme URI"

# 2. Clear the "date modified" values in column Y for the listed rows
$yRows = @(1,2,5,8,9,18,20,21,22,23,24,25,26,27,28,29,30,32,33,35,37,38,39,40,41)
foreach ($r in $yRows) {
    $ws.Cells.Item($r, 25).Value = ""
}

# 3. Clear the stray "test" value in B46
$ws.Range("B46").Value = ""
